# Weekly update: insert the newest week's 4 price rows (Abate Fettel / Winter
# Nelis, Primera / Segunda) at the top of the detail block (row 104), pushing
# the rest of the historical rows down by four rows (104-230 -> 108-234).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 104 (shifts existing rows 104..230 down
# to 108..234, and grows the sheet dimension to A1:T234 automatically).
$ws.Rows.Item(104).Resize(4).Insert()

# Common (constant) field values shared by every data row in this block.
$mercadoId = 11
$mercado = "Vega Monumental Concepción"
$region = "Bíobío"
$codreg = 8
$tipo = "Fruta"
$productoId = 100104
$producto = "Frutos de pepita"
$categoriaId = 100104005
$categoria = "Pera"
$unidad = "`$/caja 16 kilos empedrada"
$origen = "Región de O'Higgins"

# New rows: Fecha, Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, PrecioKg, KgUnidad
$newRows = @(
    @(44447, "Abate Fettel", "Primera", 50, 9000, 9000, 9000, 562, 16),
    @(44447, "Abate Fettel", "Segunda", 50, 8000, 8000, 8000, 500, 16),
    @(44447, "Winter Nelis", "Primera", 50, 9000, 9000, 9000, 562, 16),
    @(44447, "Winter Nelis", "Segunda", 50, 8000, 8000, 8000, 500, 16)
)

$r = 104
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $row[0]
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $tipo
    $ws.Cells.Item($r, 7).Value = $productoId
    $ws.Cells.Item($r, 8).Value = $producto
    $ws.Cells.Item($r, 9).Value = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $row[1]
    $ws.Cells.Item($r, 12).Value = $row[2]
    $ws.Cells.Item($r, 13).Value = $row[3]
    $ws.Cells.Item($r, 14).Value = $row[4]
    $ws.Cells.Item($r, 15).Value = $row[5]
    $ws.Cells.Item($r, 16).Value = $row[6]
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $row[7]
    $ws.Cells.Item($r, 20).Value = $row[8]
    $r = $r + 1
}
